# Clustream-MicroCluster-KDD98-Throughput-xClusterSize.xlsx
# "adjust xticks, yticks and size of figures, add normalized throughput"
#
# Sheet1 holds a throughput table (rows: cluster-size 25..150, columns: p=1..p=32).
# The commit normalizes the raw throughput numbers (µs) down to a comparable
# scale by dividing every data value in B2:G7 by 1000, and nudges the saved
# selection/scroll position of the sheet (used while re-plotting the charts).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Normalize throughput values (B2:G7): new = old / 1000 -----------------
$cols = @("B", "C", "D", "E", "F", "G")
for ($r = 2; $r -le 7; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Range("$c$r")
        $old = $cell.Value()
        $cell.Value = $old / 1000
    }
}

# --- Update the saved view: scroll so column B is left-most, select C10 ----
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("C10").Select()
